# "Hjemme passive tweaks lichtwark deleted values"
# Update header row (B1:E1) and re-measured / corrected passive values
# in rows 2-3 (columns B-E) for sheet "Ark1", including clearing a few
# cells whose measurements were removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON): B2, D2, E2 deleted; C2 replaced with new value
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -2.9033381510991703
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 (STR): B3:E3 replaced with new values
$ws.Range("B3").Value = -3.105531684919832
$ws.Range("C3").Value = -0.12958465626425664
$ws.Range("D3").Value = -6.2146140962329639
$ws.Range("E3").Value = 10.743274369729917

# Selection now only spans the edited block B1:E3 (was B1:AY3)
$ws.Range("B1:E3").Select()
